$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '35.407.66'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +1.72%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = "'" + '1.889.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +0.87%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('E4').Value = "'" + '  -0.04%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('B5').Value = "'" + 'XRP'
$ws.Range('B5').Style = 'Normal'
$ws.Range('C5').Value = "'" + 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('C5').Style = 'Normal'
$ws.Range('D5').Value = "'" + '0.693'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +0.68%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('B6').Value = "'" + 'BNB'
$ws.Range('B6').Style = 'Normal'
$ws.Range('C6').Value = "'" + 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('C6').Style = 'Normal'
$ws.Range('D6').Value = "'" + '246.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.26%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('E7').Value = "'" + '  +0.00%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = "'" + '43.22'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +2.72%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = "'" + '0.356'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +2.75%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').Value = "'" + '54.85'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +7.69%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = "'" + '0.0746'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +1.23%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('E12').Value = "'" + '  +1.92%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = "'" + '13.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +7.79%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = "'" + '2.160.80'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +0.74%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = "'" + '0.772'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +8.07%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').Value = "'" + '5.03'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +2.83%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').Value = "'" + '1.865.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.43%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = "'" + '35.441.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +1.87%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').Value = "'" + '73.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +1.10%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = "'" + '0.0₃0827'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +0.97%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = "'" + '245.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +0.64%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').Value = "'" + '12.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.96%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('D23').Value = "'" + '5.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +4.69%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = "'" + '2.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +8.38%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('E25').Value = "'" + '  -0.09%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('E26').Value = "'" + '  -4.00%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').Value = "'" + '165.94'
$ws.Range('D27').Style = 'Normal'

# Row 28
$ws.Range('D28').Value = "'" + '8.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +3.11%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('D29').Value = "'" + '18.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +0.64%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('E30').Value = "'" + '  +0.47%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('D31').Value = "'" + '0.0597'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +3.74%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').Value = "'" + '4.30'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +0.71%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = "'" + '1.89'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +23.38%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = "'" + '4.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +1.20%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('E35').Value = "'" + '  +0.03%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('E36').Value = "'" + '  -13.82%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('D37').Value = "'" + '0.858'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +3.82%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('E38').Value = "'" + '  -1.89%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = "'" + '0.0725'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +10.31%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').Value = "'" + '0.0221'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +5.38%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = "'" + '17.28'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +2.52%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = "'" + '98.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +0.56%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').Value = "'" + '1.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -0.22%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = "'" + '13.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +14.44%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = "'" + '1.327.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +3.53%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = "'" + '2.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +3.07%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').Value = "'" + '0.0809'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +2.21%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('E48').Value = "'" + '  +0.36%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('E49').Value = "'" + '  +0.61%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').Value = "'" + '6.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -1.98%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').Value = "'" + '2.063.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.53%  '
$ws.Range('E51').Style = 'Normal'
